# TC03_Verify_BLP_Solutions_ContactUS.xlsx - "Changes for New UI Prod"
#
# The test case steps for QuickOrder / ServicesMenu / SolutionsMenu /
# ResourcesMenu "verify element present" rows are removed from the
# TC03_Verify_BLP_Sol_Contact_etc sheet (rows 3, 5, 6, 7 in the original
# layout), which shifts the remaining CLICK / VERIFY_PAGE_TITLE /
# VERIFY_TEXT_PRESENT rows up so the sheet ends with 9 used rows instead
# of 13. The Testdata sheet itself is untouched data-wise; only its view
# selection/active-tab state changes (TC03 sheet becomes the active tab).

$wb = $excel.ActiveWorkbook
$wsMain = $wb.Worksheets.Item(1)      # TC03_Verify_BLP_Sol_Contact_etc
$wsData = $wb.Worksheets.Item(2)      # Testdata

# Remove the four obsolete verification-step rows. Deleting from the
# bottom up keeps the remaining row numbers stable while we work.
$wsMain.Rows(7).Delete()   # ResourcesMenu verify-present step
$wsMain.Rows(6).Delete()   # SolutionsMenu verify-present step
$wsMain.Rows(5).Delete()   # ServicesMenu verify-present step
$wsMain.Rows(3).Delete()   # QuickOrder verify-present step

# Update the selection/active-cell shown on each sheet and make the
# first sheet the active tab (it was the Testdata sheet before).
$wsData.Range("B5").Select() | Out-Null
$wsMain.Activate() | Out-Null
$wsMain.Range("C6").Select() | Out-Null
